$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text edits (volume number + report week range) ---
$ws.Range("A8").Value = "Volume 30   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/30/2023  Through  11/5/2023"

# --- Plain value updates (no type/style change) ---
$ws.Range("L14").Value = -16.666666666666
$ws.Range("N14").Value = -58.333333333333
$ws.Range("H15").Value = -100
$ws.Range("J15").Value = 15
$ws.Range("K15").Value = 40
$ws.Range("M15").Value = 16.666666666666
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 20
$ws.Range("F16").Value = 27
$ws.Range("G16").Value = 31
$ws.Range("H16").Value = -12.903225806451
$ws.Range("I16").Value = 246
$ws.Range("J16").Value = 227
$ws.Range("K16").Value = 8.370044052863
$ws.Range("L16").Value = 33.695652173913
$ws.Range("M16").Value = -6.106870229007
$ws.Range("N16").Value = -58.234295415959
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -40
$ws.Range("F17").Value = 33
$ws.Range("G17").Value = 31
$ws.Range("H17").Value = 6.451612903225
$ws.Range("I17").Value = 353
$ws.Range("J17").Value = 314
$ws.Range("K17").Value = 12.420382165605
$ws.Range("L17").Value = 42.914979757085
$ws.Range("M17").Value = 60.454545454545
$ws.Range("N17").Value = 22.569444444444
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 11.111111111111
$ws.Range("I18").Value = 207
$ws.Range("J18").Value = 113
$ws.Range("K18").Value = 83.185840707964
$ws.Range("L18").Value = 55.63909774436
$ws.Range("M18").Value = -36.697247706422
$ws.Range("N18").Value = -83.360128617363
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 114.285714285714
$ws.Range("G19").Value = 50
$ws.Range("H19").Value = 18
$ws.Range("I19").Value = 573
$ws.Range("J19").Value = 522
$ws.Range("K19").Value = 9.770114942528
$ws.Range("L19").Value = 55.70652173913
$ws.Range("M19").Value = 49.608355091383
$ws.Range("N19").Value = 12.795275590551
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 11
$ws.Range("E20").Value = -18.181818181818
$ws.Range("F20").Value = 39
$ws.Range("G20").Value = 37
$ws.Range("H20").Value = 5.405405405405
$ws.Range("I20").Value = 442
$ws.Range("J20").Value = 300
$ws.Range("K20").Value = 47.333333333333
$ws.Range("L20").Value = 99.099099099099
$ws.Range("M20").Value = 126.666666666667
$ws.Range("N20").Value = -72.8
$ws.Range("C21").Value = 39
$ws.Range("D21").Value = 35
$ws.Range("E21").Value = 11.428571428571
$ws.Range("F21").Value = 168
$ws.Range("G21").Value = 160
$ws.Range("H21").Value = 5
$ws.Range("I21").Value = 1847
$ws.Range("J21").Value = 1498
$ws.Range("K21").Value = 23.297730307076
$ws.Range("L21").Value = 56.658184902459
$ws.Range("M21").Value = 30.530035335689
$ws.Range("N21").Value = -56.996507566938
$ws.Range("J22").Value = 19
$ws.Range("K22").Value = -42.105263157894
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 200
$ws.Range("F23").Value = 9
$ws.Range("H23").Value = 50
$ws.Range("I23").Value = 104
$ws.Range("J23").Value = 94
$ws.Range("K23").Value = 10.63829787234
$ws.Range("L23").Value = 35.064935064935
$ws.Range("M23").Value = 85.714285714285
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 37
$ws.Range("E24").Value = -18.918918918918
$ws.Range("F24").Value = 146
$ws.Range("G24").Value = 127
$ws.Range("H24").Value = 14.960629921259
$ws.Range("I24").Value = 1383
$ws.Range("J24").Value = 1138
$ws.Range("K24").Value = 21.52899824253
$ws.Range("L24").Value = 61.188811188811
$ws.Range("M24").Value = 73.743718592964
$ws.Range("C25").Value = 6
$ws.Range("E25").Value = -45.454545454545
$ws.Range("F25").Value = 30
$ws.Range("H25").Value = -11.764705882352
$ws.Range("I25").Value = 446
$ws.Range("J25").Value = 455
$ws.Range("K25").Value = -1.978021978021
$ws.Range("L25").Value = 30.02915451895
$ws.Range("M25").Value = -17.560073937153
$ws.Range("F26").Value = 1
$ws.Range("H26").Value = -80
$ws.Range("J26").Value = 34
$ws.Range("K26").Value = -8.823529411764
$ws.Range("L26").Value = 3.333333333333
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -66.666666666666
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = -88.888888888888
$ws.Range("I27").Value = 60
$ws.Range("J27").Value = 47
$ws.Range("K27").Value = 27.659574468085
$ws.Range("L27").Value = 27.659574468085
$ws.Range("L28").Value = -23.529411764705
$ws.Range("L29").Value = -7.692307692307
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = -50
$ws.Range("I30").Value = 4
$ws.Range("J30").Value = 3
$ws.Range("K30").Value = 33.333333333333

# --- Text -> Number conversions (explicit number format to match sibling column style) ---
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E15").Value = -100
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 1
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E22").Value = -100
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 1
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("D30").Value = 1
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E30").Value = -100
$ws.Range("F30").NumberFormat = "#,##0"
$ws.Range("F30").Value = 1

# --- Number -> Text conversions (force text format, set literal value, then restore cosmetic style) ---
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "0"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"

# Fix cosmetic style (right-aligned General/text style) via Copy + PasteSpecial(formats only)
$ws.Range("C15").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("F28").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("F29").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").PasteSpecial(-4122)
$excel.CutCopyMode = 0
